# Update TPM-derived edge/specificity values per updated script run (see commit: "update scripts wuth new tpm")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.1332293333333333
$ws.Range("H2").Value = 0.399688
$ws.Range("I2").Value = 0.2494559176252732
$ws.Range("J2").Value = 0.2494559176252732
$ws.Range("M2").Value = 0.3007906666666667
$ws.Range("N2").Value = 0.902372
$ws.Range("O2").Value = 0.03537029821880876
$ws.Range("P2").Value = 0.03537029821880876
$ws.Range("Q2").Value = 0.04007413999288889
$ws.Range("R2").Value = 0.360667259936
$ws.Range("S2").Value = 0.008823330198852505
$ws.Range("T2").Value = 0.008823330198852505

# Row 3
$ws.Range("G3").Value = 0.1332293333333333
$ws.Range("H3").Value = 0.399688
$ws.Range("I3").Value = 0.2494559176252732
$ws.Range("J3").Value = 0.2494559176252732
$ws.Range("O3").Value = 0.899334434508434
$ws.Range("P3").Value = 0.899334434508434
$ws.Range("Q3").Value = 1.018935543205333
$ws.Range("R3").Value = 9.170419888848
$ws.Range("S3").Value = 0.2243442966123075
$ws.Range("T3").Value = 0.2243442966123075

# Row 4
$ws.Range("G4").Value = 0.1332293333333333
$ws.Range("H4").Value = 0.399688
$ws.Range("I4").Value = 0.2494559176252732
$ws.Range("J4").Value = 0.2494559176252732
$ws.Range("M4").Value = 0.0008990000000000001
$ws.Range("N4").Value = 0.002697
$ws.Range("O4").Value = 0.0001057143775473167
$ws.Range("P4").Value = 0.0001057143775473167
$ws.Range("Q4").Value = 0.0001197731706666667
$ws.Range("R4").Value = 0.001077958536
$ws.Range("S4").Value = 0.00002637107705725046
$ws.Range("T4").Value = 0.00002637107705725046

# Row 5
$ws.Range("G5").Value = 0.1332293333333333
$ws.Range("H5").Value = 0.399688
$ws.Range("I5").Value = 0.2494559176252732
$ws.Range("J5").Value = 0.2494559176252732
$ws.Range("M5").Value = 0.5528646666666667
$ws.Range("N5").Value = 1.658594
$ws.Range("O5").Value = 0.06501195117304938
$ws.Range("P5").Value = 0.06501195117304936
$ws.Range("Q5").Value = 0.07365779096355556
$ws.Range("R5").Value = 0.6629201186720001
$ws.Range("S5").Value = 0.01621761593648249
$ws.Range("T5").Value = 0.01621761593648248

# Row 6
$ws.Range("G6").Value = 0.1332293333333333
$ws.Range("H6").Value = 0.399688
$ws.Range("I6").Value = 0.2494559176252732
$ws.Range("J6").Value = 0.2494559176252732
$ws.Range("M6").Value = 0.001510333333333333
$ws.Range("N6").Value = 0.004531
$ws.Range("O6").Value = 0.0001776017221605087
$ws.Range("P6").Value = 0.0001776017221605087
$ws.Range("Q6").Value = 0.0002012207031111111
$ws.Range("R6").Value = 0.001810986328
$ws.Range("S6").Value = 0.0000443038005733785
$ws.Range("T6").Value = 0.00004430380057337851

# Row 7
$ws.Range("I7").Value = 0.4545982216136294
$ws.Range("J7").Value = 0.4545982216136294
$ws.Range("M7").Value = 0.3007906666666667
$ws.Range("N7").Value = 0.902372
$ws.Range("O7").Value = 0.03537029821880876
$ws.Range("P7").Value = 0.03537029821880876
$ws.Range("Q7").Value = 0.07302946727777777
$ws.Range("R7").Value = 0.6572652055
$ws.Range("S7").Value = 0.01607927466821419
$ws.Range("T7").Value = 0.01607927466821419

# Row 8
$ws.Range("I8").Value = 0.4545982216136294
$ws.Range("J8").Value = 0.4545982216136294
$ws.Range("O8").Value = 0.899334434508434
$ws.Range("P8").Value = 0.899334434508434
$ws.Range("S8").Value = 0.4088358345634332
$ws.Range("T8").Value = 0.4088358345634332

# Row 9
$ws.Range("I9").Value = 0.4545982216136294
$ws.Range("J9").Value = 0.4545982216136294
$ws.Range("M9").Value = 0.0008990000000000001
$ws.Range("N9").Value = 0.002697
$ws.Range("O9").Value = 0.0001057143775473167
$ws.Range("P9").Value = 0.0001057143775473167
$ws.Range("Q9").Value = 0.0002182697083333333
$ws.Range("R9").Value = 0.001964427375
$ws.Range("S9").Value = 0.00004805756803200197
$ws.Range("T9").Value = 0.00004805756803200197

# Row 10
$ws.Range("I10").Value = 0.4545982216136294
$ws.Range("J10").Value = 0.4545982216136294
$ws.Range("M10").Value = 0.5528646666666667
$ws.Range("N10").Value = 1.658594
$ws.Range("O10").Value = 0.06501195117304938
$ws.Range("P10").Value = 0.06501195117304936
$ws.Range("Q10").Value = 0.1342309338611111
$ws.Range("R10").Value = 1.20807840475
$ws.Range("S10").Value = 0.02955431738690036
$ws.Range("T10").Value = 0.02955431738690035

# Row 11
$ws.Range("I11").Value = 0.4545982216136294
$ws.Range("J11").Value = 0.4545982216136294
$ws.Range("M11").Value = 0.001510333333333333
$ws.Range("N11").Value = 0.004531
$ws.Range("O11").Value = 0.0001776017221605087
$ws.Range("P11").Value = 0.0001776017221605087
$ws.Range("Q11").Value = 0.0003666963472222222
$ws.Range("R11").Value = 0.003300267125
$ws.Range("S11").Value = 0.00008073742704968516
$ws.Range("T11").Value = 0.00008073742704968517

# Row 12
$ws.Range("G12").Value = 0.07627099999999999
$ws.Range("H12").Value = 0.228813
$ws.Range("I12").Value = 0.1428082826594534
$ws.Range("J12").Value = 0.1428082826594534
$ws.Range("M12").Value = 0.3007906666666667
$ws.Range("N12").Value = 0.902372
$ws.Range("O12").Value = 0.03537029821880876
$ws.Range("P12").Value = 0.03537029821880876
$ws.Range("Q12").Value = 0.02294160493733333
$ws.Range("R12").Value = 0.206474444436
$ws.Range("S12").Value = 0.005051171545780803
$ws.Range("T12").Value = 0.005051171545780803

# Row 13
$ws.Range("G13").Value = 0.07627099999999999
$ws.Range("H13").Value = 0.228813
$ws.Range("I13").Value = 0.1428082826594534
$ws.Range("J13").Value = 0.1428082826594534
$ws.Range("O13").Value = 0.899334434508434
$ws.Range("P13").Value = 0.899334434508434
$ws.Range("Q13").Value = 0.5833192351219999
$ws.Range("R13").Value = 5.249873116098
$ws.Range("S13").Value = 0.1284324061286601
$ws.Range("T13").Value = 0.1284324061286601

# Row 14
$ws.Range("G14").Value = 0.07627099999999999
$ws.Range("H14").Value = 0.228813
$ws.Range("I14").Value = 0.1428082826594534
$ws.Range("J14").Value = 0.1428082826594534
$ws.Range("M14").Value = 0.0008990000000000001
$ws.Range("N14").Value = 0.002697
$ws.Range("O14").Value = 0.0001057143775473167
$ws.Range("P14").Value = 0.0001057143775473167
$ws.Range("Q14").Value = 0.000068567629
$ws.Range("R14").Value = 0.000617108661
$ws.Range("S14").Value = 0.00001509688870994538
$ws.Range("T14").Value = 0.00001509688870994538

# Row 15
$ws.Range("G15").Value = 0.07627099999999999
$ws.Range("H15").Value = 0.228813
$ws.Range("I15").Value = 0.1428082826594534
$ws.Range("J15").Value = 0.1428082826594534
$ws.Range("M15").Value = 0.5528646666666667
$ws.Range("N15").Value = 1.658594
$ws.Range("O15").Value = 0.06501195117304938
$ws.Range("P15").Value = 0.06501195117304936
$ws.Range("Q15").Value = 0.04216754099133333
$ws.Range("R15").Value = 0.3795078689219999
$ws.Range("S15").Value = 0.00928424509936342
$ws.Range("T15").Value = 0.009284245099363419

# Row 16
$ws.Range("G16").Value = 0.07627099999999999
$ws.Range("H16").Value = 0.228813
$ws.Range("I16").Value = 0.1428082826594534
$ws.Range("J16").Value = 0.1428082826594534
$ws.Range("M16").Value = 0.001510333333333333
$ws.Range("N16").Value = 0.004531
$ws.Range("O16").Value = 0.0001776017221605087
$ws.Range("P16").Value = 0.0001776017221605087
$ws.Range("Q16").Value = 0.0001151946336666667
$ws.Range("R16").Value = 0.001036751703
$ws.Range("S16").Value = 0.00002536299693910363
$ws.Range("T16").Value = 0.00002536299693910364

# Row 17
$ws.Range("G17").Value = 0.08178766666666666
$ws.Range("H17").Value = 0.245363
$ws.Range("I17").Value = 0.153137578101644
$ws.Range("J17").Value = 0.153137578101644
$ws.Range("M17").Value = 0.3007906666666667
$ws.Range("N17").Value = 0.902372
$ws.Range("O17").Value = 0.03537029821880876
$ws.Range("P17").Value = 0.03537029821880876
$ws.Range("Q17").Value = 0.02460096678177778
$ws.Range("R17").Value = 0.221408701036
$ws.Range("S17").Value = 0.005416521805961268
$ws.Range("T17").Value = 0.005416521805961268

# Row 18
$ws.Range("G18").Value = 0.08178766666666666
$ws.Range("H18").Value = 0.245363
$ws.Range("I18").Value = 0.153137578101644
$ws.Range("J18").Value = 0.153137578101644
$ws.Range("O18").Value = 0.899334434508434
$ws.Range("P18").Value = 0.899334434508434
$ws.Range("Q18").Value = 0.6255106024886666
$ws.Range("R18").Value = 5.629595422398
$ws.Range("S18").Value = 0.1377218972040332
$ws.Range("T18").Value = 0.1377218972040332

# Row 19
$ws.Range("G19").Value = 0.08178766666666666
$ws.Range("H19").Value = 0.245363
$ws.Range("I19").Value = 0.153137578101644
$ws.Range("J19").Value = 0.153137578101644
$ws.Range("M19").Value = 0.0008990000000000001
$ws.Range("N19").Value = 0.002697
$ws.Range("O19").Value = 0.0001057143775473167
$ws.Range("P19").Value = 0.0001057143775473167
$ws.Range("Q19").Value = 0.00007352711233333334
$ws.Range("R19").Value = 0.0006617440110000001
$ws.Range("S19").Value = 0.00001618884374811889
$ws.Range("T19").Value = 0.00001618884374811889

# Row 20
$ws.Range("G20").Value = 0.08178766666666666
$ws.Range("H20").Value = 0.245363
$ws.Range("I20").Value = 0.153137578101644
$ws.Range("J20").Value = 0.153137578101644
$ws.Range("M20").Value = 0.5528646666666667
$ws.Range("N20").Value = 1.658594
$ws.Range("O20").Value = 0.06501195117304938
$ws.Range("P20").Value = 0.06501195117304936
$ws.Range("Q20").Value = 0.04521751106911111
$ws.Range("R20").Value = 0.406957599622
$ws.Range("S20").Value = 0.009955772750303116
$ws.Range("T20").Value = 0.009955772750303115

# Row 21
$ws.Range("G21").Value = 0.08178766666666666
$ws.Range("H21").Value = 0.245363
$ws.Range("I21").Value = 0.153137578101644
$ws.Range("J21").Value = 0.153137578101644
$ws.Range("M21").Value = 0.001510333333333333
$ws.Range("N21").Value = 0.004531
$ws.Range("O21").Value = 0.0001776017221605087
$ws.Range("P21").Value = 0.0001776017221605087
$ws.Range("Q21").Value = 0.0001235266392222222
$ws.Range("R21").Value = 0.001111739753
$ws.Range("S21").Value = 0.00002719749759834138
$ws.Range("T21").Value = 0.00002719749759834138
